$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1273.0435
$ws.Range("I17").Value = 499
$ws.Range("J17").Value = 1308.2273
$ws.Range("K17").Value = 1497
$ws.Range("L17").Value = 3924.6819
$ws.Range("N17").Value = -4260.6819
$ws.Range("M17").Value = -1329
$ws.Range("H51").Value = 250004500
$ws.Range("I51").Value = 250004500
$ws.Range("K51").Value = 250004500
$ws.Range("M51").Value = -250004016
$ws.Range("H55").Value = 467.4
$ws.Range("J55").Value = 313.57144
$ws.Range("L55").Value = 313.57144
$ws.Range("N55").Value = -741.5714399999999
$ws.Range("H69").Value = 7999
$ws.Range("J69").Value = 7999
$ws.Range("L69").Value = 23997
$ws.Range("N69").Value = -25745
$ws.Range("H72").Value = 7999
$ws.Range("J72").Value = 7999
$ws.Range("L72").Value = 71991
$ws.Range("N72").Value = -80727
$ws.Range("H111").Value = 1874.2858
$ws.Range("I111").Value = 7120
$ws.Range("K111").Value = 21360
$ws.Range("M111").Value = -18293
$ws.Range("H138").Value = 2460.457
$ws.Range("J138").Value = 2644.049
$ws.Range("L138").Value = 7932.147
$ws.Range("N138").Value = -18212.147

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4162.3335
$ws.Range("I32").Value = 4212.9565
$ws.Range("K32").Value = 4212.9565
$ws.Range("M32").Value = -3925.9565
$ws.Range("H61").Value = 3320.5
$ws.Range("I61").Value = 1908.25
$ws.Range("J61").Value = 3885.4
$ws.Range("K61").Value = 1908.25
$ws.Range("L61").Value = 3885.4
$ws.Range("M61").Value = -1696.25
$ws.Range("N61").Value = -4309.4
$ws.Range("H136").Value = 3320.5
$ws.Range("I136").Value = 1908.25
$ws.Range("J136").Value = 3885.4
$ws.Range("K136").Value = 5724.75
$ws.Range("L136").Value = 11656.2
$ws.Range("M136").Value = -3174.75
$ws.Range("N136").Value = -16756.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 607.913
$ws.Range("I22").Value = 582
$ws.Range("J22").Value = 636.1818
$ws.Range("K22").Value = 582
$ws.Range("L22").Value = 636.1818
$ws.Range("M22").Value = -409
$ws.Range("N22").Value = -982.1818
$ws.Range("H86").Value = 2460.7083
$ws.Range("I86").Value = 2298.0908
$ws.Range("J86").Value = 4249.5
$ws.Range("K86").Value = 2298.0908
$ws.Range("L86").Value = 4249.5
$ws.Range("M86").Value = -1175.0908
$ws.Range("N86").Value = -6495.5
$ws.Range("H89").Value = 2460.7083
$ws.Range("I89").Value = 2298.0908
$ws.Range("J89").Value = 4249.5
$ws.Range("K89").Value = 11490.454
$ws.Range("L89").Value = 21247.5
$ws.Range("M89").Value = -5874.454
$ws.Range("N89").Value = -32479.5
$ws.Range("H99").Value = 80357.766
$ws.Range("I99").Value = 94104.73
$ws.Range("K99").Value = 94104.73
$ws.Range("M99").Value = -92606.73

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H31").Value = 4152.383
$ws.Range("I31").Value = 2980.6667
$ws.Range("K31").Value = 2980.6667
$ws.Range("M31").Value = -2685.6667
$ws.Range("H34").Value = 4152.383
$ws.Range("I34").Value = 2980.6667
$ws.Range("K34").Value = 2980.6667
$ws.Range("M34").Value = -2778.6667
$ws.Range("H58").Value = 2464.2942
$ws.Range("I58").Value = 1625.6364
$ws.Range("K58").Value = 1625.6364
$ws.Range("M58").Value = -1422.6364
$ws.Range("H107").Value = 7143958.5
$ws.Range("I107").Value = 12500578
$ws.Range("J107").Value = 1799
$ws.Range("K107").Value = 12500578
$ws.Range("L107").Value = 1799
$ws.Range("M107").Value = -12498658
$ws.Range("N107").Value = -5639
$ws.Range("H132").Value = 12350542
$ws.Range("I132").Value = 4698.8184
$ws.Range("K132").Value = 14096.4552
$ws.Range("M132").Value = -11566.4552
$ws.Range("H133").Value = 99999
$ws.Range("J133").Value = 99999
$ws.Range("L133").Value = 99999
$ws.Range("N133").Value = -105059
$ws.Range("H136").Value = 2464.2942
$ws.Range("I136").Value = 1625.6364
$ws.Range("K136").Value = 4876.9092
$ws.Range("M136").Value = -2326.9092

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 6803.4287
$ws.Range("I141").Value = 6803.4287
$ws.Range("K141").Value = 20410.2861
$ws.Range("M141").Value = -15230.2861

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 38466960
$ws.Range("I70").Value = 100003800
$ws.Range("K70").Value = 100003800
$ws.Range("M70").Value = -100003530
$ws.Range("H73").Value = 38466960
$ws.Range("I73").Value = 100003800
$ws.Range("K73").Value = 100003800
$ws.Range("M73").Value = -100002864
$ws.Range("H126").Value = 9548.25
$ws.Range("I126").Value = 2224.5
$ws.Range("K126").Value = 6673.5
$ws.Range("M126").Value = -4203.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2118.5454
$ws.Range("I7").Value = 1811.5555
$ws.Range("K7").Value = 1811.5555
$ws.Range("M7").Value = -1699.5555
$ws.Range("H55").Value = 459.8125
$ws.Range("I55").Value = 338.3158
$ws.Range("J55").Value = 637.38464
$ws.Range("K55").Value = 338.3158
$ws.Range("L55").Value = 637.38464
$ws.Range("M55").Value = -165.3158
$ws.Range("N55").Value = -983.38464
$ws.Range("H82").Value = 856.6842
$ws.Range("I82").Value = 795
$ws.Range("J82").Value = 990.3333
$ws.Range("K82").Value = 795
$ws.Range("L82").Value = 990.3333
$ws.Range("M82").Value = -434
$ws.Range("N82").Value = -1712.3333
$ws.Range("H85").Value = 856.6842
$ws.Range("I85").Value = 795
$ws.Range("J85").Value = 990.3333
$ws.Range("K85").Value = 795
$ws.Range("L85").Value = 990.3333
$ws.Range("M85").Value = 453
$ws.Range("N85").Value = -3486.3333
$ws.Range("H126").Value = 2118.5454
$ws.Range("I126").Value = 1811.5555
$ws.Range("K126").Value = 5434.666499999999
$ws.Range("M126").Value = -2964.666499999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H81").Value = 4055
$ws.Range("I81").Value = 3513.0625
$ws.Range("J81").Value = 5500.1665
$ws.Range("K81").Value = 7026.125
$ws.Range("L81").Value = 11000.333
$ws.Range("M81").Value = -5965.125
$ws.Range("N81").Value = -13122.333
$ws.Range("H84").Value = 4055
$ws.Range("I84").Value = 3513.0625
$ws.Range("J84").Value = 5500.1665
$ws.Range("K84").Value = 35130.625
$ws.Range("L84").Value = 55001.665
$ws.Range("M84").Value = -29826.625
$ws.Range("N84").Value = -65609.66500000001
$ws.Range("H100").Value = 90910210
$ws.Range("I100").Value = 1213.5
$ws.Range("K100").Value = 2427
$ws.Range("M100").Value = -1886
$ws.Range("H126").Value = 2934.5454
$ws.Range("I126").Value = 3060
$ws.Range("K126").Value = 9180
$ws.Range("M126").Value = -6710
